$wb = $excel.ActiveWorkbook

# --- PortalPage (sheet1): update the demo login credentials shown in B2/C2 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "kikosana@calliduscloud.com"
$ws1.Range("C2").Value = "Kiran09@@"

# --- Commissions Data (sheet2): populate headers (row1) + sample values (row2) ---
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A1").Value = "Credit Rule"
$ws2.Range("B1").Value = "Rule Type"
$ws2.Range("C1").Value = "Deposit Rule"
$ws2.Range("D1").Value = "Earning Code"
$ws2.Range("E1").Value = "Earning Type"
$ws2.Range("G1").Value = "Fixed Value Type"
$ws2.Range("I1").Value = "Incentive Rule"
$ws2.Range("K1").Value = "Measurement Rule"
$ws2.Range("L1").Value = "Plan Components"
$ws2.Range("M1").Value = "Plans"
$ws2.Range("N1").Value = "Position Group"
$ws2.Range("O1").Value = "Rate Table"
$ws2.Range("P1").Value = "Rate Table variable"
$ws2.Range("Q1").Value = "Reason Codes"
$ws2.Range("R1").Value = "Territory Elements"
$ws2.Range("S1").Value = "Territory Variables"
$ws2.Range("T1").Value = "Unit Type"

$ws2.Range("A2").Value = "CR Territory Booking"
$ws2.Range("B2").Value = "Bonus123"
$ws2.Range("C2").Value = "DR Base Commissions Flat Rate MTD"
$ws2.Range("E2").Value = "commission"

$ws2.Range("F1").Value = "Event Type"
$ws2.Range("F2").Value = "invoicing"

$ws2.Range("G2").Value = "Fix1"

$ws2.Range("H1").Value = "Fixed value Variable"
$ws2.Range("H2").Value = "Var_FV_Demo1"

$ws2.Range("I2").Value = "IR Base Commissions Flat Rate Quota Based MTD"

$ws2.Range("J1").Value = "Lookup Table Variable"
$ws2.Range("J2").Value = "Var_LT_Demo4"

$ws2.Range("K2").Value = "MR Testrule"
$ws2.Range("L2").Value = "Aggregated Revenue against Territory based on Flat Rate"
$ws2.Range("M2").Value = "compensation plan"
$ws2.Range("N2").Value = "P_G_Demo1"
$ws2.Range("O2").Value = "R_T_Demo1"
$ws2.Range("P2").Value = "Var_RT_Demo2"
$ws2.Range("Q2").Value = "Test RC"
$ws2.Range("R2").Value = "Territory_Demo"
$ws2.Range("S2").Value = "Var_T_Demo3"
$ws2.Range("T2").Value = "USD"

# Column widths to fit the new data.
# NOTE: the host snaps ColumnWidth to its internal pixel grid (~1/6-character
# steps), the same way Excel quantizes a width typed into the Format >
# Column Width dialog to whole pixels. That rounding adds a small constant
# offset, so the requested widths are pre-compensated (-5/6) to land on the
# desired stored width after the host's own rounding is applied.
$colWidthCorrection = 5/6
$ws2.Columns.Item(1).ColumnWidth = 24 - $colWidthCorrection
$ws2.Columns.Item(2).ColumnWidth = 19.28515625 - $colWidthCorrection
$ws2.Columns.Item(3).ColumnWidth = 36.85546875 - $colWidthCorrection
$ws2.Range($ws2.Columns.Item(4), $ws2.Columns.Item(8)).ColumnWidth = 22.140625 - $colWidthCorrection
$ws2.Columns.Item(9).ColumnWidth = 43.7109375 - $colWidthCorrection
$ws2.Columns.Item(10).ColumnWidth = 23.7109375 - $colWidthCorrection
$ws2.Columns.Item(11).ColumnWidth = 18.140625 - $colWidthCorrection
$ws2.Columns.Item(12).ColumnWidth = 52.7109375 - $colWidthCorrection
$ws2.Range($ws2.Columns.Item(13), $ws2.Columns.Item(20)).ColumnWidth = 18.140625 - $colWidthCorrection

# Rename the worksheet
$ws2.Name = "Commissions Data"

# Restore prior selections before switching the active tab
$ws1.Range("H16").Select()
$ws2.Range("F13:F14").Select()

# Commissions Data is now the active/visible tab
$ws2.Activate()
